# Moodle endpoints workbook update: add "Activities" and "Grades" rows to
# the Endpoint sheet, tidy up the existing error message for the "get
# enrolled users" endpoint, and switch the active sheet/selection to match.

$wb = $excel.ActiveWorkbook
$wsGeneral = $wb.Worksheets.Item("General")
$wsEndpoint = $wb.Worksheets.Item("Endpoint")

# --- Endpoint sheet: expand the old single-cell "Roles" / "Actividades" /
# "Calificaciones" rows (19-21) into a richer table describing the
# activities & grades endpoints (rows 19-22).
# Values are written in an order that mirrors how the new shared strings
# end up appended to the workbook's string table.

# Row 20: Actividades
$wsEndpoint.Range("B20").Value = "Activities"

# Row 21 & 22: Calificaciones (Grades)
$wsEndpoint.Range("A22").Value = "Calificaciones"
$wsEndpoint.Range("B21").Value = "Grades"
$wsEndpoint.Range("B22").Value = "Grades"
$wsEndpoint.Range("D21").Value = "OBTENER GRADES POR CURSO"
$wsEndpoint.Range("C21").Value = "Permite obtener calificaciones de un curso"
$wsEndpoint.Range("D22").Value = "OBTENER GRADES POR USUARIO"
$wsEndpoint.Range("C22").Value = "Permite obtener calificaciones finales por usuario"
$wsEndpoint.Range("E21").Value = "gradereport_user_get_grade_items"
$wsEndpoint.Range("E22").Value = "gradereport_overview_get_course_grades"
$wsEndpoint.Range("F21").Value = "POST"
$wsEndpoint.Range("G21").Value = "application/x-www-form-urlencoded"
$wsEndpoint.Range("H21").Value = "courseid=int"
$wsEndpoint.Range("I21").Value = "Objeto usergrades que contiene un array de grades"
$wsEndpoint.Range("J21").Value = "No"
$wsEndpoint.Range("K21").Value = "courseid"
$wsEndpoint.Range("L21").Value = "Presenta error específico de curso no existente"
$wsEndpoint.Range("M21").Value = "Presenta información de todos los estudiantes tengan o no calificaciones"
$wsEndpoint.Range("K22").Value = "userid"
$wsEndpoint.Range("H22").Value = "userid=int"
$wsEndpoint.Range("F22").Value = "POST"
$wsEndpoint.Range("G22").Value = "application/x-www-form-urlencoded"
$wsEndpoint.Range("I22").Value = "Objeto grades que contiene un array de courseid, grade, rawgrade y un objeto warnings"
$wsEndpoint.Range("J22").Value = "No"
$wsEndpoint.Range("L22").Value = "Presenta error específico de user no existente"

# Row 19: Roles (unchanged entity label, new columns)
$wsEndpoint.Range("B19").Value = "Role"
$wsEndpoint.Range("M20").Value = "No existe un endpoint para listar las actividades pero se puede utilizar el endpoint gradereport_user_get_grade_items y filtrar los datos de las actividades"
$wsEndpoint.Range("M19").Value = "No existe un endpoint para listar los roles existentes"

# --- Endpoint sheet: replace the generic "course doesn't exist" error text
# on the "get enrolled users" row with a more specific message. (This also
# drops the old shared string "Presenta error si no existe el curso" since
# it becomes unused.)
$wsEndpoint.Range("L17").Value = "Presenta error específico de curso inválido"

# --- Switch the active sheet/selection to the Endpoint tab.
$wsEndpoint.Activate()
$wsEndpoint.Range("F23").Select()
